$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 57.8
$ws.Range("N2").Value = 85.8724807945396

$ws.Range("K3").Value = 53
$ws.Range("N3").Value = 85.8724807945396

$ws.Range("D4").Value = 90923.57000000001
$ws.Range("K4").Value = 50.8
$ws.Range("N4").Value = 85.8724807945396

$ws.Range("K5").Value = 50.2
$ws.Range("N5").Value = 85.8724807945396

$ws.Range("K6").Value = 47.8
$ws.Range("N6").Value = 85.8724807945396
